$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.121.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.124.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.78%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.46%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.118.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.94%  "

$ws.Range("E9").Value = "  -1.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.10%  "

$ws.Range("E12").Value = "  -2.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000250"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.635.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.95%  "

$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.093.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.125.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("E19").Value = "  -2.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.53%  "

$ws.Range("E22").Value = "  -3.82%  "

$ws.Range("E23").Value = "  -2.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.73%  "

$ws.Range("E27").Value = "  -3.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.88%  "

$ws.Range("E30").Value = "  -1.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -9.40%  "

$ws.Range("E34").Value = "  -4.39%  "

$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("E36").Value = "  -1.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0713"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0389"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "421.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.55%  "

$ws.Range("E41").Value = "  -0.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("E43").Value = "  -12.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.862.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  -3.81%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.94%  "

$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.91%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.114"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.40%  "

